# Populate PÇ (program outcome) average support data for rows 6-8 on the
# "EK5-PÇ Karşılama Yüzdeleri" sheet so the AVERAGEIF() formulas in row 19
# (and the bar3D chart that reads from it) no longer resolve to #DIV/0!.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EK5-PÇ Karşılama Yüzdeleri")

$rows = 6, 7, 8

foreach ($r in $rows) {
    $ws.Cells.Item($r, 1).Value = "test"
    $ws.Cells.Item($r, 2).Value = "test2"

    # Block E:L (skip M)
    foreach ($c in "E", "F", "G", "H", "I", "K", "L") {
        $ws.Range($c + $r).Value = 0.9248000383377075
    }
    $ws.Range("J" + $r).Value = 0.9371429085731506

    # Block P:W (skip N, O)
    foreach ($c in "P", "Q", "R", "S", "T", "U", "V", "W") {
        $ws.Range($c + $r).Value = 0.6489999890327454
    }

    # Block AA:AH (skip Y, Z)
    foreach ($c in "AA", "AB", "AC", "AD", "AE", "AF", "AG", "AH") {
        $ws.Range($c + $r).Value = 0.6800000071525574
    }

    # Block AL:AS, excluding AQ (skip AJ, AK)
    foreach ($c in "AL", "AM", "AN", "AO", "AP", "AR", "AS") {
        $ws.Range($c + $r).Value = 0.6840000152587891
    }

    # Block AW:BD (skip AU, AV)
    foreach ($c in "AW", "AX", "AY", "AZ", "BA", "BC", "BD") {
        $ws.Range($c + $r).Value = 0.7344500422477722
    }
    $ws.Range("BB" + $r).Value = 0.7553809285163879
}

Write-Host "Cells done"
Write-Host $ws.Shapes.Count

